# Auto-generated Excel COM-interop script
# Applies the numeric updates captured in the commit diff for Sheets/Golem_Profits.xlsx
# (workbook tabs: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)

$wb = $excel.ActiveWorkbook


# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 153.66667
$ws.Range("I6").Value = 169.375
$ws.Range("J6").Value = 28
$ws.Range("K6").Value = 508.125
$ws.Range("L6").Value = 84
$ws.Range("M6").Value = -396.125
$ws.Range("N6").Value = -308
# Row 8
$ws.Range("H8").Value = 10
$ws.Range("I8").Value = 10
$ws.Range("K8").Value = 30
$ws.Range("M8").Value = 109
# Row 15
$ws.Range("H15").Value = 1957.6364
$ws.Range("I15").Value = 1957.6364
$ws.Range("K15").Value = 5872.9092
$ws.Range("M15").Value = -5703.9092
# Row 31
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()
# Row 39
$ws.Range("H39").Value = 8022.933
$ws.Range("I39").Value = 1918.125
$ws.Range("J39").Value = 14999.857
$ws.Range("K39").Value = 5754.375
$ws.Range("L39").Value = 44999.571
$ws.Range("M39").Value = -5458.375
$ws.Range("N39").Value = -45591.571
# Row 64
$ws.Range("H64").Value = 7083.3335
$ws.Range("J64").Value = 10000
$ws.Range("L64").Value = 10000
$ws.Range("N64").Value = -10496
# Row 67
$ws.Range("H67").Value = 7083.3335
$ws.Range("J67").Value = 10000
$ws.Range("L67").Value = 10000
$ws.Range("N67").Value = -11716
# Row 70
$ws.Range("H70").Value = 3666.5833
$ws.Range("I70").Value = 3388.889
$ws.Range("J70").Value = 4499.6665
$ws.Range("K70").Value = 10166.667
$ws.Range("L70").Value = 13498.9995
$ws.Range("M70").Value = -9896.667000000001
$ws.Range("N70").Value = -14038.9995
# Row 73
$ws.Range("H73").Value = 3666.5833
$ws.Range("I73").Value = 3388.889
$ws.Range("J73").Value = 4499.6665
$ws.Range("K73").Value = 10166.667
$ws.Range("L73").Value = 13498.9995
$ws.Range("M73").Value = -9230.667000000001
$ws.Range("N73").Value = -15370.9995
# Row 107
$ws.Range("H107").Value = 75968
$ws.Range("I107").Value = 100936.11
$ws.Range("K107").Value = 100936.11
$ws.Range("M107").Value = -99016.11
# Row 137
$ws.Range("H137").Value = 939.8
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
# Row 138
$ws.Range("H138").Value = 5566.3447
$ws.Range("I138").Value = 2228.4285
$ws.Range("K138").Value = 6685.2855
$ws.Range("M138").Value = -1545.2855

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 3
$ws.Range("H3").Value = 666.3333
$ws.Range("I3").Value = 666.3333
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 666.3333
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -551.3333
$ws.Range("N3").ClearContents()
# Row 10
$ws.Range("H10").Value = 2749.5
$ws.Range("I10").Value = 1999.6666
$ws.Range("J10").Value = 4999
$ws.Range("K10").Value = 1999.6666
$ws.Range("L10").Value = 4999
$ws.Range("M10").Value = -1829.6666
$ws.Range("N10").Value = -5339
# Row 12
$ws.Range("H12").Value = 24121.6
$ws.Range("J12").Value = 40000
$ws.Range("L12").Value = 40000
$ws.Range("N12").Value = -40346
# Row 21
$ws.Range("H21").Value = 2937.5
$ws.Range("J21").Value = 2937.5
$ws.Range("L21").Value = 2937.5
$ws.Range("N21").Value = -3685.5
# Row 25
$ws.Range("H25").Value = 10069.857
$ws.Range("I25").Value = 4000
$ws.Range("J25").Value = 11081.5
$ws.Range("K25").Value = 4000
$ws.Range("L25").Value = 11081.5
$ws.Range("M25").Value = -3598
$ws.Range("N25").Value = -11885.5
# Row 27
$ws.Range("H27").Value = 8500
$ws.Range("J27").Value = 8500
$ws.Range("L27").Value = 8500
$ws.Range("N27").Value = -8868
# Row 30
$ws.Range("H30").Value = 14440
$ws.Range("J30").Value = 17425
$ws.Range("L30").Value = 17425
$ws.Range("N30").Value = -17725
# Row 35
$ws.Range("H35").Value = 2783.6
$ws.Range("I35").Value = 2783.6
$ws.Range("K35").Value = 2783.6
$ws.Range("M35").Value = -2377.6
# Row 45
$ws.Range("H45").Value = 2833
$ws.Range("I45").Value = 3249.5
$ws.Range("K45").Value = 3249.5
$ws.Range("M45").Value = -2872.5
# Row 92
$ws.Range("H92").Value = 225000
$ws.Range("J92").Value = 225000
$ws.Range("L92").Value = 225000
$ws.Range("N92").Value = -229992
# Row 128
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 2146.7273
$ws.Range("I99").Value = 2146.7273
$ws.Range("K99").Value = 2146.7273
$ws.Range("M99").Value = -648.7273

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 2
$ws.Range("H2").Value = 919.2143
$ws.Range("I2").Value = 685.625
$ws.Range("J2").Value = 1230.6666
$ws.Range("K2").Value = 685.625
$ws.Range("L2").Value = 1230.6666
$ws.Range("M2").Value = -572.625
$ws.Range("N2").Value = -1456.6666
# Row 13
$ws.Range("H13").Value = 475
$ws.Range("I13").Value = 450
$ws.Range("J13").Value = 500
$ws.Range("K13").Value = 450
$ws.Range("L13").Value = 500
$ws.Range("M13").Value = -311
$ws.Range("N13").Value = -778
# Row 63
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
# Row 66
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
# Row 132
$ws.Range("H132").Value = 1333.3334
$ws.Range("I132").Value = 1000
$ws.Range("J132").Value = 1500
$ws.Range("K132").Value = 3000
$ws.Range("L132").Value = 4500
$ws.Range("M132").Value = -470
$ws.Range("N132").Value = -9560

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 1556.3265
$ws.Range("I4").Value = 380.18518
$ws.Range("J4").Value = 2999.7727
$ws.Range("K4").Value = 1140.55554
$ws.Range("L4").Value = 8999.3181
$ws.Range("M4").Value = -1028.55554
$ws.Range("N4").Value = -9223.3181
# Row 11
$ws.Range("H11").Value = 333.33334
$ws.Range("I11").Value = 375
$ws.Range("J11").Value = 250
$ws.Range("K11").Value = 1125
$ws.Range("L11").Value = 750
$ws.Range("M11").Value = -985
$ws.Range("N11").Value = -1030

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 57
$ws.Range("H57").Value = 10703.333
$ws.Range("I57").Value = 1055
$ws.Range("J57").Value = 30000
$ws.Range("K57").Value = 1055
$ws.Range("L57").Value = 30000
$ws.Range("M57").Value = -235
$ws.Range("N57").Value = -31640
# Row 80
$ws.Range("H80").Value = 4125
$ws.Range("I80").Value = 4166.6665
$ws.Range("K80").Value = 4166.6665
$ws.Range("M80").Value = -3168.6665
# Row 83
$ws.Range("H83").Value = 4125
$ws.Range("I83").Value = 4166.6665
$ws.Range("K83").Value = 20833.3325
$ws.Range("M83").Value = -15841.3325
# Row 132
$ws.Range("H132").Value = 240
$ws.Range("I132").Value = 240
$ws.Range("K132").Value = 720
$ws.Range("M132").Value = 1810

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 39899.8
$ws.Range("I2").Value = 750
$ws.Range("J2").Value = 65999.664
$ws.Range("K2").Value = 750
$ws.Range("L2").Value = 65999.664
$ws.Range("M2").Value = -638
$ws.Range("N2").Value = -66223.664
# Row 26
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
# Row 46
$ws.Range("H46").Value = 500
$ws.Range("J46").Value = 500
$ws.Range("L46").Value = 500
$ws.Range("N46").Value = -876
# Row 55
$ws.Range("H55").Value = 1815.2
$ws.Range("I55").Value = 1901
$ws.Range("J55").Value = 1793.75
$ws.Range("K55").Value = 1901
$ws.Range("L55").Value = 1793.75
$ws.Range("M55").Value = -1728
$ws.Range("N55").Value = -2139.75
# Row 132
$ws.Range("H132").Value = 3653.75
$ws.Range("I132").Value = 3653.75
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 10961.25
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -8431.25
$ws.Range("N132").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
# Row 4
$ws.Range("H4").Value = 4201.5
$ws.Range("J4").Value = 4201.5
$ws.Range("L4").Value = 4201.5
$ws.Range("N4").Value = -4427.5
